# Add two new columns, "I0" (column I) and "IF" (column J), to the stats
# sheet -- mirrors the existing header style (bold + border, from H1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from H1 onto the two new header cells, then
# set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-37: (I0, IF) pairs per row.
$data = @{
    2  = @(1, 3)
    3  = @(1, 2)
    4  = @(1, 2)
    5  = @(1, 3)
    6  = @(1, 4)
    7  = @(8, 8)
    8  = @(7, 9)
    9  = @(9, 9)
    10 = @(5, 6)
    11 = @(6, 6)
    12 = @(6, 7)
    13 = @(6, 7)
    14 = @(8, 8)
    15 = @(2, 2)
    16 = @(6, 6)
    17 = @(9, 9)
    18 = @(6, 7)
    19 = @(5, 6)
    20 = @(6, 7)
    21 = @(5, 5)
    22 = @(5, 6)
    23 = @(3, 4)
    24 = @(6, 8)
    25 = @(6, 6)
    26 = @(6, 7)
    27 = @(8, 8)
    28 = @(6, 7)
    29 = @(8, 9)
    30 = @(10, 10)
    31 = @(5, 7)
    32 = @(1, 4)
    33 = @(5, 5)
    34 = @(1, 2)
    35 = @(6, 6)
    36 = @(4, 5)
    37 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
